$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells for the prereq/coreq columns.
$ws.Range("E1").Value = "Prerequisites"
$ws.Range("F1").Value = "CoRequisites"

# Reflect the header-row selection left behind when the workbook was saved.
$ws.Range("A1:I1").Select()
